# Auto-generated edit script applying the diff to before.xlsx
$wb = $excel.ActiveWorkbook

$sheetData = @{
    'runs' = @(
        ,@('B2', 8534.034)
        ,@('C2', 5831.925000000001)
        ,@('D2', 97)
        ,@('E2', 8620.1532)
        ,@('F2', 154.4010000000003)
        ,@('H2', 5954.944461062575)
        ,@('J2', 1134634)
        ,@('B3', 8515.2375)
        ,@('C3', 5988.5885)
        ,@('D3', 92)
        ,@('E3', 8592.07105)
        ,@('F3', 120)
        ,@('H3', 5988.5885)
        ,@('I3', 59)
        ,@('J3', 1133912)
        ,@('B4', 8504.9935)
        ,@('C4', 5980.991)
        ,@('D4', 87)
        ,@('E4', 8582.50165)
        ,@('F4', 126)
        ,@('H4', 5980.991)
        ,@('J4', 1133693)
        ,@('B5', 8515.2375)
        ,@('C5', 5980.991)
        ,@('D5', 92)
        ,@('E5', 8592.07105)
        ,@('F5', 126)
        ,@('H5', 5980.991)
        ,@('J5', 1133912)
    );
    'per_resource_all_runs' = @(
        ,@('E6', 1796)
        ,@('G6', 0.000522300148)
        ,@('H6', 0.0006880141944000001)
        ,@('E7', 24970)
        ,@('G7', 0.007261600609999998)
        ,@('H7', 0.009565542558000002)
        ,@('E10', 1189)
        ,@('G10', 0.000345776657)
        ,@('H10', 0.0004554837846)
        ,@('E14', 2370)
        ,@('G14', 0.0006892268099999998)
        ,@('H14', 0.0009079029180000002)
        ,@('E15', 16312)
        ,@('G15', 0.004743741655999999)
        ,@('H15', 0.0062488237968)
        ,@('E16', 3298)
        ,@('G16', 0.0009591012739999998)
        ,@('H16', 0.0012634024572000001)
        ,@('E29', 6373)
        ,@('F29', 14809)
        ,@('G29', 0.001853351249)
        ,@('H29', 0.0024413777622)
        ,@('E30', 3430)
        ,@('G30', 0.0009974885899999997)
        ,@('H30', 0.0013139692020000002)
        ,@('E31', 39093)
        ,@('F31', 131282)
        ,@('G31', 0.011368752608999999)
        ,@('H31', 0.014975801170200004)
        ,@('E32', 24064)
        ,@('F32', 278993)
        ,@('G32', 0.006998124031999999)
        ,@('H32', 0.009218470809600003)
        ,@('B33', 'https://res.cloudinary.com/ssp/image/fetch/c_fill,w_360,ar_16:9/c_limit,w_750/https://aimblob.blob.core.windows.net/aimfiles/f0463a60-75f8-4054-a627-698ba68e0507.jpg')
        ,@('E33', 30424)
        ,@('F33', 29846)
        ,@('G33', 0.008847694711999999)
        ,@('H33', 0.0116548685136)
        ,@('B34', 'https://res.cloudinary.com/ssp/image/fetch/c_fill,w_360,ar_16:9/c_limit,w_750/https://aimblob.blob.core.windows.net/aimfiles/c6fdd324-eb71-446d-94ae-ccc28f251fbe.jpg')
        ,@('E34', 19709)
        ,@('F34', 19133)
        ,@('G34', 0.005731633416999999)
        ,@('H34', 0.007550151312600001)
        ,@('B35', 'https://res.cloudinary.com/ssp/image/fetch/c_fill,w_360,ar_16:9/c_limit,w_750/https://aimblob.blob.core.windows.net/aimfiles/07c4945c-6349-4d51-ae3d-9488eb7fe7f1.jpg')
        ,@('D35', 'image/jpeg')
        ,@('E35', 12500)
        ,@('F35', 11933)
        ,@('G35', 0.0036351624999999992)
        ,@('H35', 0.0047885175000000006)
        ,@('E39', 468)
        ,@('G39', 0.00013610048399999996)
        ,@('H39', 0.00017928209520000002)
        ,@('E40', 1027)
        ,@('G40', 0.00029866495099999997)
        ,@('H40', 0.00039342459780000006)
        ,@('E46', 1198)
        ,@('G46', 0.00034839397399999997)
        ,@('H46', 0.0004589315172000001)
        ,@('E48', 924)
        ,@('G48', 0.00026871121199999995)
        ,@('H48', 0.00035396721360000005)
        ,@('E50', 2253)
        ,@('G50', 0.0006552016889999999)
        ,@('H50', 0.0008630823942000002)
        ,@('E51', 16240)
        ,@('G51', 0.00472280312)
        ,@('H51', 0.006221241936)
        ,@('E52', 3172)
        ,@('G52', 0.0009224588359999998)
        ,@('H52', 0.0012151342008000002)
        ,@('E60', 2025)
        ,@('G60', 0.0005888963249999999)
        ,@('H60', 0.000775739835)
        ,@('E65', 6310)
        ,@('F65', 14809)
        ,@('G65', 0.0018350300299999997)
        ,@('H65', 0.0024172436340000004)
        ,@('E66', 3403)
        ,@('G66', 0.0009896366389999997)
        ,@('H66', 0.0013036260042)
        ,@('E67', 38895)
        ,@('F67', 131282)
        ,@('G67', 0.011311171634999999)
        ,@('H67', 0.014899951053)
        ,@('E68', 24064)
        ,@('F68', 278993)
        ,@('G68', 0.006998124031999999)
        ,@('H68', 0.009218470809600003)
        ,@('B69', 'https://res.cloudinary.com/ssp/image/fetch/c_fill,w_360,ar_16:9/c_limit,w_750/https://aimblob.blob.core.windows.net/aimfiles/f0463a60-75f8-4054-a627-698ba68e0507.jpg')
        ,@('E69', 30424)
        ,@('F69', 29846)
        ,@('G69', 0.008847694711999999)
        ,@('H69', 0.0116548685136)
        ,@('B70', 'https://res.cloudinary.com/ssp/image/fetch/c_fill,w_360,ar_16:9/c_limit,w_750/https://aimblob.blob.core.windows.net/aimfiles/c6fdd324-eb71-446d-94ae-ccc28f251fbe.jpg')
        ,@('E70', 19709)
        ,@('F70', 19133)
        ,@('G70', 0.005731633416999999)
        ,@('H70', 0.007550151312600001)
        ,@('B71', 'https://res.cloudinary.com/ssp/image/fetch/c_fill,w_360,ar_16:9/c_limit,w_750/https://aimblob.blob.core.windows.net/aimfiles/07c4945c-6349-4d51-ae3d-9488eb7fe7f1.jpg')
        ,@('D71', 'image/jpeg')
        ,@('E71', 12499)
        ,@('F71', 11933)
        ,@('G71', 0.0036348716869999994)
        ,@('H71', 0.004788134418600001)
        ,@('E75', 839)
        ,@('G75', 0.00024399210699999996)
        ,@('H75', 0.00032140529459999997)
        ,@('E76', 656)
        ,@('G76', 0.00019077332799999997)
        ,@('H76', 0.00025130139840000005)
        ,@('E79', 24970)
        ,@('G79', 0.007261600609999998)
        ,@('H79', 0.009565542558000002)
        ,@('E82', 1063)
        ,@('G82', 0.00030913421899999996)
        ,@('H82', 0.00040721552820000007)
        ,@('E84', 925)
        ,@('G84', 0.00026900202499999996)
        ,@('H84', 0.00035435029500000003)
        ,@('E85', 1963)
        ,@('G85', 0.000570865919)
        ,@('H85', 0.0007519887882000002)
        ,@('E86', 2262)
        ,@('G86', 0.0006578190059999999)
        ,@('H86', 0.0008665301268000001)
        ,@('E87', 16285)
        ,@('G87', 0.004735889704999999)
        ,@('H87', 0.006238480599000001)
        ,@('E88', 3190)
        ,@('G88', 0.0009276934699999998)
        ,@('H88', 0.0012220296659999999)
        ,@('E96', 1890)
        ,@('G96', 0.0005496365699999999)
        ,@('H96', 0.000724023846)
        ,@('E101', 6260)
        ,@('F101', 14809)
        ,@('G101', 0.0018204893799999998)
        ,@('H101', 0.002398089564)
        ,@('E102', 3304)
        ,@('G102', 0.0009608461519999999)
        ,@('H102', 0.0012657009456000002)
        ,@('E103', 38976)
        ,@('F103', 131282)
        ,@('G103', 0.011334727487999999)
        ,@('H103', 0.014930980646400001)
        ,@('E104', 24064)
        ,@('F104', 278993)
        ,@('G104', 0.006998124031999999)
        ,@('H104', 0.009218470809600003)
        ,@('B105', 'https://res.cloudinary.com/ssp/image/fetch/c_fill,w_360,ar_16:9/c_limit,w_750/https://aimblob.blob.core.windows.net/aimfiles/f0463a60-75f8-4054-a627-698ba68e0507.jpg')
        ,@('E105', 30424)
        ,@('F105', 29846)
        ,@('G105', 0.008847694711999999)
        ,@('H105', 0.0116548685136)
        ,@('B106', 'https://res.cloudinary.com/ssp/image/fetch/c_fill,w_360,ar_16:9/c_limit,w_750/https://aimblob.blob.core.windows.net/aimfiles/c6fdd324-eb71-446d-94ae-ccc28f251fbe.jpg')
        ,@('E106', 19709)
        ,@('F106', 19133)
        ,@('G106', 0.005731633416999999)
        ,@('H106', 0.007550151312600001)
        ,@('B107', 'https://res.cloudinary.com/ssp/image/fetch/c_fill,w_360,ar_16:9/c_limit,w_750/https://aimblob.blob.core.windows.net/aimfiles/07c4945c-6349-4d51-ae3d-9488eb7fe7f1.jpg')
        ,@('D107', 'image/jpeg')
        ,@('E107', 12499)
        ,@('F107', 11933)
        ,@('G107', 0.0036348716869999994)
        ,@('H107', 0.004788134418600001)
    );
    'summary_by_type' = @(
        ,@('C3', 311857)
        ,@('C4', 1296621)
        ,@('C5', 282480)
        ,@('D5', 1547580)
        ,@('C6', 16545)
        ,@('C7', 763372)
        ,@('D7', 745578)
    );
    'co2' = @(
        ,@('B2', 1134634)
        ,@('C2', 0.32996631744199995)
        ,@('D2', 0.43465718120760005)
        ,@('B3', 1133912)
        ,@('C3', 0.329756350456)
        ,@('D3', 0.43438059643680005)
        ,@('B4', 1133693)
        ,@('C4', 0.32969266240899997)
        ,@('D4', 0.43429670161020006)
        ,@('B5', 1133912)
        ,@('C5', 0.329756350456)
        ,@('D5', 0.43438059643680005)
    )
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($pair in $sheetData[$sheetName]) {
        $ws.Range($pair[0]).Value = $pair[1]
    }
}

Write-Output "Applied $($sheetData.Values | ForEach-Object { $_.Count } | Measure-Object -Sum | Select-Object -ExpandProperty Sum) cell edits"